$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-10-15 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-16 Thursday", 2) | Out-Null

# Update the arithmetic table cells (row-major order) via Tables(1).Cell(row, col)
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "98-19="
$t.Cell(1, 2).Range.Text = "56+26="
$t.Cell(1, 3).Range.Text = "23+8="
$t.Cell(1, 4).Range.Text = "27+27="
$t.Cell(1, 5).Range.Text = "91-34="
$t.Cell(2, 1).Range.Text = "90-62="
$t.Cell(2, 2).Range.Text = "22-19="
$t.Cell(2, 3).Range.Text = "9+2="
$t.Cell(2, 4).Range.Text = "29+48="
$t.Cell(2, 5).Range.Text = "90-79="
$t.Cell(3, 1).Range.Text = "81-44="
$t.Cell(3, 2).Range.Text = "63-25="
$t.Cell(3, 3).Range.Text = "13+59="
$t.Cell(3, 4).Range.Text = "41-17="
$t.Cell(3, 5).Range.Text = "7+49="
$t.Cell(4, 1).Range.Text = "80-16="
$t.Cell(4, 2).Range.Text = "74-37="
$t.Cell(4, 3).Range.Text = "34+27="
$t.Cell(4, 4).Range.Text = "34-17="
$t.Cell(4, 5).Range.Text = "28+68="
$t.Cell(5, 1).Range.Text = "3+89="
$t.Cell(5, 2).Range.Text = "23+38="
$t.Cell(5, 3).Range.Text = "69+7="
$t.Cell(5, 4).Range.Text = "81-25="
$t.Cell(5, 5).Range.Text = "37+9="
$t.Cell(6, 1).Range.Text = "42-35="
$t.Cell(6, 2).Range.Text = "70-24="
$t.Cell(6, 3).Range.Text = "71-59="
$t.Cell(6, 4).Range.Text = "39+57="
$t.Cell(6, 5).Range.Text = "66+27="
$t.Cell(7, 1).Range.Text = "71-33="
$t.Cell(7, 2).Range.Text = "9+67="
$t.Cell(7, 4).Range.Text = "18+63="
$t.Cell(7, 5).Range.Text = "59+2="
$t.Cell(8, 1).Range.Text = "26+55="
$t.Cell(8, 2).Range.Text = "54+29="
$t.Cell(8, 3).Range.Text = "53-45="
$t.Cell(8, 4).Range.Text = "27-8="
$t.Cell(8, 5).Range.Text = "41-16="
$t.Cell(9, 1).Range.Text = "18+28="
$t.Cell(9, 2).Range.Text = "9+39="
$t.Cell(9, 3).Range.Text = "19+77="
$t.Cell(9, 4).Range.Text = "49+22="
$t.Cell(9, 5).Range.Text = "94-58="
$t.Cell(10, 1).Range.Text = "27+66="
$t.Cell(10, 2).Range.Text = "80-29="
$t.Cell(10, 3).Range.Text = "63+28="
$t.Cell(10, 4).Range.Text = "17+38="
$t.Cell(10, 5).Range.Text = "16+47="
$t.Cell(11, 1).Range.Text = "61-2="
$t.Cell(11, 2).Range.Text = "78+9="
$t.Cell(11, 3).Range.Text = "86-69="
$t.Cell(11, 4).Range.Text = "58-29="
$t.Cell(11, 5).Range.Text = "18+45="
$t.Cell(12, 1).Range.Text = "71-22="
$t.Cell(12, 2).Range.Text = "77+8="
$t.Cell(12, 3).Range.Text = "94-56="
$t.Cell(12, 4).Range.Text = "80-77="
$t.Cell(12, 5).Range.Text = "14+38="
$t.Cell(13, 1).Range.Text = "66+29="
$t.Cell(13, 2).Range.Text = "18+48="
$t.Cell(13, 3).Range.Text = "40-33="
$t.Cell(13, 4).Range.Text = "82-66="
$t.Cell(13, 5).Range.Text = "46+29="
$t.Cell(14, 1).Range.Text = "29+57="
$t.Cell(14, 2).Range.Text = "86-37="
$t.Cell(14, 3).Range.Text = "51-28="
$t.Cell(14, 4).Range.Text = "39+28="
$t.Cell(14, 5).Range.Text = "56+28="
$t.Cell(15, 1).Range.Text = "22-19="
$t.Cell(15, 2).Range.Text = "37+35="
$t.Cell(15, 3).Range.Text = "67+4="
$t.Cell(15, 4).Range.Text = "29+46="
$t.Cell(15, 5).Range.Text = "90-36="
$t.Cell(16, 1).Range.Text = "48-19="
$t.Cell(16, 2).Range.Text = "57+25="
$t.Cell(16, 3).Range.Text = "81-2="
$t.Cell(16, 4).Range.Text = "98-59="
$t.Cell(16, 5).Range.Text = "27+55="
$t.Cell(17, 1).Range.Text = "64-47="
$t.Cell(17, 2).Range.Text = "49+37="
$t.Cell(17, 3).Range.Text = "58+4="
$t.Cell(17, 4).Range.Text = "41-15="
$t.Cell(17, 5).Range.Text = "70-27="
$t.Cell(18, 1).Range.Text = "48+16="
$t.Cell(18, 2).Range.Text = "26+5="
$t.Cell(18, 3).Range.Text = "68-59="
$t.Cell(18, 4).Range.Text = "17+47="
$t.Cell(18, 5).Range.Text = "57-39="
$t.Cell(19, 1).Range.Text = "28+43="
$t.Cell(19, 2).Range.Text = "47-38="
$t.Cell(19, 3).Range.Text = "91-38="
$t.Cell(19, 4).Range.Text = "73-58="
$t.Cell(19, 5).Range.Text = "58+18="
$t.Cell(20, 1).Range.Text = "26-18="
$t.Cell(20, 2).Range.Text = "70-2="
$t.Cell(20, 3).Range.Text = "59+38="
$t.Cell(20, 4).Range.Text = "64+27="
$t.Cell(20, 5).Range.Text = "63-37="
